$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.933.84'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '1.899.63'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7582'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.72'
$ws.Range('E6').Value = '  -1.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3054'
$ws.Range('E8').Value = '  -2.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.49'
$ws.Range('E9').Value = '  -6.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06845'
$ws.Range('E10').Value = '  -2.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07998'
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7514'
$ws.Range('E12').Value = '  -4.13%  '
$ws.Range('D13').Value = '1.898.92'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.229'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.23'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').Value = '29.941.10'
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.98'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.969'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.90'
$ws.Range('E19').Value = '  -1.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007736'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.150.73'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.974'
$ws.Range('E24').Value = '  +4.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.264'
$ws.Range('E25').Value = '  -2.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.74'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1306'
$ws.Range('E28').Value = '  +2.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.026'
$ws.Range('E29').Value = '  -4.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.383'
$ws.Range('E30').Value = '  +2.43%  '
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.299'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.038'
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05362'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.259'
$ws.Range('E35').Value = '  -4.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7287'
$ws.Range('E36').Value = '  -3.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.728'
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01927'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.780'
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('E40').Value = '  -3.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4418'
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.43'
$ws.Range('E42').Value = '  -4.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.914'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8280'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.11'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.593'
$ws.Range('E47').Value = '  -2.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.805'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').Value = '2.060.82'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.27'
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05961'
$ws.Range('E51').Value = '  -0.65%  '
